# Auto-generated edit script: refresh market-price derived columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 10652.333
$ws.Range("I33").Value = 1630.8462
$ws.Range("J33").Value = 25312.25
$ws.Range("K33").Value = 1630.8462
$ws.Range("L33").Value = 25312.25
$ws.Range("M33").Value = -1401.8462
$ws.Range("N33").Value = -25770.25
# Row 131
$ws.Range("H131").Value = 1584.75
$ws.Range("I131").Value = 779.6667
$ws.Range("K131").Value = 2339.0001
$ws.Range("M131").Value = 2700.9999
# Row 138
$ws.Range("H138").Value = 3326.6836
$ws.Range("I138").Value = 697.6786
$ws.Range("J138").Value = 6832.024
$ws.Range("K138").Value = 2093.0358
$ws.Range("L138").Value = 20496.072
$ws.Range("M138").Value = 3046.9642
$ws.Range("N138").Value = -30776.072

$ws = $wb.Worksheets.Item("ARM")
# Row 17
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()
# Row 61
$ws.Range("H61").Value = 1008.14636
$ws.Range("I61").Value = 908.36365
$ws.Range("J61").Value = 1419.75
$ws.Range("K61").Value = 908.36365
$ws.Range("L61").Value = 1419.75
$ws.Range("M61").Value = -696.36365
$ws.Range("N61").Value = -1843.75
# Row 74
$ws.Range("H74").Value = 1229.0476
$ws.Range("I74").Value = 1107.8572
$ws.Range("J74").Value = 1471.4286
$ws.Range("K74").Value = 1107.8572
$ws.Range("L74").Value = 1471.4286
$ws.Range("M74").Value = -233.8571999999999
$ws.Range("N74").Value = -3219.4286
# Row 77
$ws.Range("H77").Value = 1229.0476
$ws.Range("I77").Value = 1107.8572
$ws.Range("J77").Value = 1471.4286
$ws.Range("K77").Value = 5539.286
$ws.Range("L77").Value = 7357.143
$ws.Range("M77").Value = -1171.286
$ws.Range("N77").Value = -16093.143
# Row 92
$ws.Range("H92").Value = 20730
$ws.Range("J92").Value = 20730
$ws.Range("L92").Value = 20730
$ws.Range("N92").Value = -25722
# Row 101
$ws.Range("H101").Value = 33099.832
$ws.Range("J101").Value = 33099.832
$ws.Range("L101").Value = 33099.832
$ws.Range("N101").Value = -39589.832
# Row 110
$ws.Range("H110").Value = 1451.56
$ws.Range("I110").Value = 909.7646999999999
$ws.Range("J110").Value = 2602.875
$ws.Range("K110").Value = 909.7646999999999
$ws.Range("L110").Value = 2602.875
$ws.Range("M110").Value = 1135.2353
$ws.Range("N110").Value = -6692.875
# Row 132
$ws.Range("H132").Value = 2102.1333
$ws.Range("I132").Value = 1481.2222
$ws.Range("J132").Value = 3033.5
$ws.Range("K132").Value = 4443.6666
$ws.Range("L132").Value = 9100.5
$ws.Range("M132").Value = -1913.6666
$ws.Range("N132").Value = -14160.5
# Row 136
$ws.Range("H136").Value = 1008.14636
$ws.Range("I136").Value = 908.36365
$ws.Range("J136").Value = 1419.75
$ws.Range("K136").Value = 2725.09095
$ws.Range("L136").Value = 4259.25
$ws.Range("M136").Value = -175.0909499999998
$ws.Range("N136").Value = -9359.25

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1404.5
$ws.Range("I94").Value = 856.3
$ws.Range("J94").Value = 2775
$ws.Range("K94").Value = 856.3
$ws.Range("L94").Value = 2775
$ws.Range("M94").Value = -405.3
$ws.Range("N94").Value = -3677
# Row 100
$ws.Range("H100").Value = 23300
$ws.Range("J100").Value = 23300
$ws.Range("L100").Value = 23300
$ws.Range("N100").Value = -25464
# Row 134
$ws.Range("H134").Value = 1421.9269
$ws.Range("I134").Value = 978.1786
$ws.Range("J134").Value = 2377.6924
$ws.Range("K134").Value = 2934.5358
$ws.Range("L134").Value = 7133.0772
$ws.Range("M134").Value = -399.5357999999997
$ws.Range("N134").Value = -12203.0772

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 495.18182
$ws.Range("I22").Value = 341.83334
$ws.Range("J22").Value = 679.2
$ws.Range("K22").Value = 341.83334
$ws.Range("L22").Value = 679.2
$ws.Range("M22").Value = 8.166659999999979
$ws.Range("N22").Value = -1379.2
# Row 58
$ws.Range("H58").Value = 1642.909
$ws.Range("I58").Value = 1296
$ws.Range("J58").Value = 2250
$ws.Range("K58").Value = 1296
$ws.Range("L58").Value = 2250
$ws.Range("M58").Value = -1093
$ws.Range("N58").Value = -2656
# Row 74
$ws.Range("H74").Value = 12104.667
$ws.Range("J74").Value = 12104.667
$ws.Range("L74").Value = 12104.667
$ws.Range("N74").Value = -13852.667
# Row 77
$ws.Range("H77").Value = 12104.667
$ws.Range("J77").Value = 12104.667
$ws.Range("L77").Value = 36314.001
$ws.Range("N77").Value = -45050.001
# Row 88
$ws.Range("H88").Value = 27744.3
$ws.Range("J88").Value = 27744.3
$ws.Range("L88").Value = 27744.3
$ws.Range("N88").Value = -28556.3
# Row 91
$ws.Range("H91").Value = 27744.3
$ws.Range("J91").Value = 27744.3
$ws.Range("L91").Value = 27744.3
$ws.Range("N91").Value = -30552.3
# Row 92
$ws.Range("H92").Value = 16000
$ws.Range("J92").Value = 16000
$ws.Range("L92").Value = 16000
$ws.Range("N92").Value = -20992
# Row 96
$ws.Range("H96").Value = 20909.857
$ws.Range("J96").Value = 20909.857
$ws.Range("L96").Value = 20909.857
$ws.Range("N96").Value = -26401.857
# Row 106
$ws.Range("H106").Value = 26617.5
$ws.Range("J106").Value = 26617.5
$ws.Range("L106").Value = 26617.5
$ws.Range("N106").Value = -29141.5
# Row 136
$ws.Range("H136").Value = 1642.909
$ws.Range("I136").Value = 1296
$ws.Range("J136").Value = 2250
$ws.Range("K136").Value = 3888
$ws.Range("L136").Value = 6750
$ws.Range("M136").Value = -1338
$ws.Range("N136").Value = -11850

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 3112723.8
$ws.Range("I113").Value = 481354.4
$ws.Range("J113").Value = 7936901
$ws.Range("K113").Value = 1444063.2
$ws.Range("L113").Value = 23810703
$ws.Range("M113").Value = -1441893.2
$ws.Range("N113").Value = -23815043
# Row 131
$ws.Range("H131").Value = 924.4
$ws.Range("J131").Value = 965.619
$ws.Range("L131").Value = 2896.857
$ws.Range("N131").Value = -12976.857

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5022.5938
$ws.Range("I70").Value = 4574
$ws.Range("J70").Value = 5879
$ws.Range("K70").Value = 4574
$ws.Range("L70").Value = 5879
$ws.Range("M70").Value = -4304
$ws.Range("N70").Value = -6419
# Row 73
$ws.Range("H73").Value = 5022.5938
$ws.Range("I73").Value = 4574
$ws.Range("J73").Value = 5879
$ws.Range("K73").Value = 4574
$ws.Range("L73").Value = 5879
$ws.Range("M73").Value = -3638
$ws.Range("N73").Value = -7751

$ws = $wb.Worksheets.Item("WVR")
# Row 28
$ws.Range("H28").Value = 70019
$ws.Range("J28").Value = 70019
$ws.Range("L28").Value = 70019
$ws.Range("N28").Value = -70715
# Row 101
$ws.Range("H101").Value = 6397.625
$ws.Range("J101").Value = 6397.625
$ws.Range("L101").Value = 6397.625
$ws.Range("N101").Value = -12887.625
